$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.607193756757511, 1.482022928510048, 1.405793804184327, 1.374887289893309, 1.369764798536437, 1.405376350504184, 1.563904813376212, 1.879767770628007, 2.114929703103201, 2.222598029120149, 2.263469412450888, 2.254662589462896, 2.225958542359763, 2.208389479164964, 2.107907266822679, 2.046442057662318, 2.011154011647989, 1.999217263584342, 2.052978393359012, 2.23438690829488, 2.353529269854221, 2.289887465070535, 2.05002316300363, 1.793778750815306)
$Cvals = @(0.224866629977214, 0.1962453532894983, 0.1786231537928131, 0.1714294769405456, 0.1702342107012953, 0.1785261882458826, 0.2150079369851596, 0.2861782521521832, 0.338268128475363, 0.3619281753240102, 0.3708828628603555, 0.3689545218383614, 0.3626649784599181, 0.3588118279603236, 0.3367211921163573, 0.3231603208703575, 0.315357068653384, 0.3127144288771149, 0.3246042482492442, 0.3645124992089563, 0.3905666756556911, 0.376663546349107, 0.3239514702789279, 0.2669617032012184)
$Dvals = @(0.07806988087429545, 0.0708611335141569, 0.06647633230835481, 0.0646997463184249, 0.06440536054627444, 0.06645233132151418, 0.07557561716375005, 0.09380275460173948, 0.1074116859867331, 0.1136527347793788, 0.1160234855130682, 0.1155125700146158, 0.1138476285534438, 0.1128287727819242, 0.1070048450647931, 0.1034450524263377, 0.1014022758559179, 0.1007114326400966, 0.103823508487352, 0.1143364599831358, 0.1212505012770606, 0.1175563355721039, 0.1036523966928655, 0.08883454536818647)
$Evals = @(0.05847840465997844, 0.05884300634309714, 0.05908547284501253, 0.05918896008643504, 0.0592064268284842, 0.05908684955587695, 0.05860026086777115, 0.05779350425069563, 0.05729048316799901, 0.05708108699837933, 0.05700458555810872, 0.05702093734932578, 0.05707473722601186, 0.05710805479838132, 0.0573045585428158, 0.05743008261072369, 0.05750410974111198, 0.05752948824343296, 0.05741653107720524, 0.05705885912829878, 0.0568413749793244, 0.05695596175197437, 0.05742265192667073, 0.057995987551549)
$Gvals = @(1.385453889518089, 1.359119116830669, 1.344082219484989, 1.338237315730339, 1.337283790651895, 1.344002250976899, 1.376137285944424, 1.448241133163265, 1.506900225563584, 1.534851846630346, 1.545621087219587, 1.543293492539902, 1.535734126097083, 1.53112790160219, 1.505099221723952, 1.489457710546105, 1.480580271476981, 1.47759491554649, 1.491110428283605, 1.537949467467428, 1.569638386876051, 1.552626103037511, 1.490362876340441, 1.427747394510874)
$Hvals = @(1.228819075792586, 1.223653090373375, 1.221276028625397, 1.220506412326216, 1.220390616646114, 1.221264844445244, 1.226872367170358, 1.24421504258251, 1.260883445033556, 1.269331764279229, 1.27265641065128, 1.27193479381242, 1.26960276521686, 1.268190696260348, 1.260348830473419, 1.255760522055482, 1.253202882157353, 1.252350866647078, 1.256240520362468, 1.2702843261917, 1.280194439783401, 1.274837960027185, 1.256023263271288, 1.23883802044864)
$Ivals = @(1.106905566371644, 1.105316970415636, 1.10512194104713, 1.105237715402048, 1.105268703635069, 1.105122713169919, 1.106195257003002, 1.114536299139303, 1.124533005274415, 1.129935187191094, 1.13210488564151, 1.131632069075366, 1.130111197844798, 1.129195803193355, 1.124197252436616, 1.121350485465655, 1.119793475177701, 1.119280068839544, 1.121645202137202, 1.130554540090969, 1.137100651864785, 1.133540312478829, 1.121511712772275, 1.111605247845247)
$Lvals = @(0.2575035652529323, 0.2472221148328515, 0.2410625738949648, 0.2385909296116182, 0.2381828315558749, 0.2410290850512098, 0.253926616804435, 0.2804428318673473, 0.3006838464057466, 0.3100600054884808, 0.3136349359147061, 0.3128639236116868, 0.31035362755064, 0.3088191789853454, 0.3000744953108523, 0.2947531667661849, 0.2917083217611633, 0.2906801055897006, 0.2953179904694281, 0.3110902992932836, 0.3215406137447019, 0.3159500177610681, 0.2950625886530105, 0.2731370802230941)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value2 = $Bvals[$i]
    $ws.Range("C" + $row).Value2 = $Cvals[$i]
    $ws.Range("D" + $row).Value2 = $Dvals[$i]
    $ws.Range("E" + $row).Value2 = $Evals[$i]
    $ws.Range("G" + $row).Value2 = $Gvals[$i]
    $ws.Range("H" + $row).Value2 = $Hvals[$i]
    $ws.Range("I" + $row).Value2 = $Ivals[$i]
    $ws.Range("L" + $row).Value2 = $Lvals[$i]
}

$wb.Save()
